$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.774.72"
$ws.Range("E2").Value = "  -0.14%  "
$ws.Range("D3").Value = "3.465.34"
$ws.Range("E3").Value = "  -0.97%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "612.43"
$ws.Range("E5").Value = "  +1.98%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.97"
$ws.Range("E6").Value = "  -1.88%  "
$ws.Range("D7").Value = "3.458.49"
$ws.Range("E7").Value = "  -1.00%  "
$ws.Range("E8").Value = "  -1.99%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("E9").Value = "  -0.04%  "
$ws.Range("E10").Value = "  +0.69%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.11"
$ws.Range("E11").Value = "  -2.25%  "
$ws.Range("E12").Value = "  -2.18%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "44.52"
$ws.Range("E13").Value = "  -2.94%  "
$ws.Range("E14").Value = "  -0.97%  "
$ws.Range("D15").Value = "4.022.56"
$ws.Range("E15").Value = "  -1.08%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.20"
$ws.Range("E16").Value = "  -0.17%  "
$ws.Range("B17").Value = "BitcoinCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "585.06"
$ws.Range("E17").Value = "  -2.73%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.467.01"
$ws.Range("E18").Value = "  -0.87%  "
$ws.Range("D19").Value = "69.803.38"
$ws.Range("E19").Value = "  -0.18%  "
$ws.Range("E20").Value = "  +1.18%  "
$ws.Range("E21").Value = "  +0.99%  "
$ws.Range("E22").Value = "  -1.26%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.95"
$ws.Range("E23").Value = "  -1.96%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "96.03"
$ws.Range("E24").Value = "  +0.96%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "15.27"
$ws.Range("E25").Value = "  -0.94%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.64"
$ws.Range("E26").Value = "  -1.33%  "
$ws.Range("E27").Value = "  +0.06%  "
$ws.Range("E28").Value = "  -3.85%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.01"
$ws.Range("E29").Value = "  -2.21%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.65"
$ws.Range("E30").Value = "  -2.96%  "
$ws.Range("E31").Value = "  -1.92%  "
$ws.Range("E32").Value = "  -5.60%  "
$ws.Range("E33").Value = "  -2.58%  "
$ws.Range("E34").Value = "  -3.80%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "587.16"
$ws.Range("E35").Value = "  -15.68%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0480"
$ws.Range("E38").Value = "  -2.84%  "
$ws.Range("E39").Value = "  +0.44%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "56.32"
$ws.Range("E40").Value = "  -0.73%  "
$ws.Range("E41").Value = "  -0.23%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.16"
$ws.Range("E42").Value = "  -10.47%  "
$ws.Range("D43").Value = "3.252.14"
$ws.Range("E43").Value = "  -1.88%  "
$ws.Range("E44").Value = "  +2.38%  "
$ws.Range("E45").Value = "  -4.74%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.79"
$ws.Range("E46").Value = "  -3.70%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "30.87"
$ws.Range("E47").Value = "  -3.67%  "
$ws.Range("E48").Value = "  -4.55%  "
$ws.Range("E49").Value = "  -1.87%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "134.26"
$ws.Range("E50").Value = "  +1.09%  "
$ws.Range("E51").Value = "  +0.00%  "
